$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1847.1177
$ws.Range("I19").Value = 1008.375
$ws.Range("J19").Value = 2592.6667
$ws.Range("K19").Value = 1008.375
$ws.Range("L19").Value = 2592.6667
$ws.Range("M19").Value = -833.375
$ws.Range("N19").Value = -2942.6667
$ws.Range("H43").Value = 9167.263000000001
$ws.Range("J43").Value = 9176.611000000001
$ws.Range("L43").Value = 9176.611000000001
$ws.Range("N43").Value = -9314.611000000001
$ws.Range("H129").Value = 23810200
$ws.Range("I129").Value = 651
$ws.Range("J129").Value = 111111880
$ws.Range("K129").Value = 1953
$ws.Range("L129").Value = 333335640
$ws.Range("M129").Value = 3047
$ws.Range("N129").Value = -333345640
$ws.Range("H132").Value = 200578.34
$ws.Range("I132").Value = 212688.44
$ws.Range("J132").Value = 15899.25
$ws.Range("K132").Value = 638065.3200000001
$ws.Range("L132").Value = 47697.75
$ws.Range("M132").Value = -635535.3200000001
$ws.Range("N132").Value = -52757.75
$ws.Range("H137").Value = 4844.8247
$ws.Range("I137").Value = 3481.9285
$ws.Range("K137").Value = 10445.7855
$ws.Range("M137").Value = -7895.7855
$ws.Range("H138").Value = 2894.7847
$ws.Range("J138").Value = 5762.0386
$ws.Range("L138").Value = 17286.1158
$ws.Range("N138").Value = -27566.1158

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4257.1113
$ws.Range("I61").Value = 3950.4595
$ws.Range("K61").Value = 3950.4595
$ws.Range("M61").Value = -3738.4595
$ws.Range("H74").Value = 7389.3335
$ws.Range("I74").Value = 10199.833
$ws.Range("K74").Value = 10199.833
$ws.Range("M74").Value = -9325.833000000001
$ws.Range("H77").Value = 7389.3335
$ws.Range("I77").Value = 10199.833
$ws.Range("K77").Value = 50999.165
$ws.Range("M77").Value = -46631.165
$ws.Range("H80").Value = 44444
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("H83").Value = 44444
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("H97").Value = 700.6177
$ws.Range("I97").Value = 613.72
$ws.Range("J97").Value = 942
$ws.Range("K97").Value = 613.72
$ws.Range("L97").Value = 942
$ws.Range("M97").Value = -117.72
$ws.Range("N97").Value = -1934
$ws.Range("H132").Value = 812553.25
$ws.Range("I132").Value = 1121438
$ws.Range("K132").Value = 3364314
$ws.Range("M132").Value = -3361784
$ws.Range("H136").Value = 4257.1113
$ws.Range("I136").Value = 3950.4595
$ws.Range("K136").Value = 11851.3785
$ws.Range("M136").Value = -9301.378499999999
$ws.Range("M80").ClearContents()
$ws.Range("M83").ClearContents()

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 11775.583
$ws.Range("I86").Value = 7283.8823
$ws.Range("J86").Value = 22684
$ws.Range("K86").Value = 7283.8823
$ws.Range("L86").Value = 22684
$ws.Range("M86").Value = -6160.8823
$ws.Range("N86").Value = -24930
$ws.Range("H89").Value = 11775.583
$ws.Range("I89").Value = 7283.8823
$ws.Range("J89").Value = 22684
$ws.Range("K89").Value = 36419.4115
$ws.Range("L89").Value = 113420
$ws.Range("M89").Value = -30803.4115
$ws.Range("N89").Value = -124652
$ws.Range("H105").Value = 2341.4736
$ws.Range("I105").Value = 2221.5557
$ws.Range("J105").Value = 4500
$ws.Range("K105").Value = 2221.5557
$ws.Range("L105").Value = 4500
$ws.Range("M105").Value = -474.5556999999999
$ws.Range("N105").Value = -7994
$ws.Range("H107").Value = 6670880.5
$ws.Range("I107").Value = 7696336
$ws.Range("K107").Value = 7696336
$ws.Range("M107").Value = -7694416

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6601.3
$ws.Range("I31").Value = 4035
$ws.Range("J31").Value = 7242.875
$ws.Range("K31").Value = 4035
$ws.Range("L31").Value = 7242.875
$ws.Range("M31").Value = -3740
$ws.Range("N31").Value = -7832.875
$ws.Range("H34").Value = 6601.3
$ws.Range("I34").Value = 4035
$ws.Range("J34").Value = 7242.875
$ws.Range("K34").Value = 4035
$ws.Range("L34").Value = 7242.875
$ws.Range("M34").Value = -3833
$ws.Range("N34").Value = -7646.875
$ws.Range("H58").Value = 18873088
$ws.Range("I58").Value = 23259254
$ws.Range("K58").Value = 23259254
$ws.Range("M58").Value = -23259051
$ws.Range("H136").Value = 18873088
$ws.Range("I136").Value = 23259254
$ws.Range("K136").Value = 69777762
$ws.Range("M136").Value = -69775212

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 2098.2856
$ws.Range("I8").Value = 2098.2856
$ws.Range("K8").Value = 6294.8568
$ws.Range("M8").Value = -6155.8568
$ws.Range("H12").Value = 3348386.8
$ws.Range("I12").Value = 10000005
$ws.Range("J12").Value = 22577.5
$ws.Range("K12").Value = 30000015
$ws.Range("L12").Value = 67732.5
$ws.Range("M12").Value = -29999842
$ws.Range("N12").Value = -68078.5
$ws.Range("H115").Value = 8324.666999999999
$ws.Range("I115").Value = 9999.5
$ws.Range("K115").Value = 29998.5
$ws.Range("M115").Value = -28823.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 10200.542
$ws.Range("I70").Value = 5895
$ws.Range("K70").Value = 5895
$ws.Range("M70").Value = -5625
$ws.Range("H73").Value = 10200.542
$ws.Range("I73").Value = 5895
$ws.Range("K73").Value = 5895
$ws.Range("M73").Value = -4959
$ws.Range("H80").Value = 43482172
$ws.Range("I80").Value = 76924850
$ws.Range("J80").Value = 6698.3
$ws.Range("K80").Value = 76924850
$ws.Range("L80").Value = 6698.3
$ws.Range("M80").Value = -76923852
$ws.Range("N80").Value = -8694.299999999999
$ws.Range("H83").Value = 43482172
$ws.Range("I83").Value = 76924850
$ws.Range("J83").Value = 6698.3
$ws.Range("K83").Value = 384624250
$ws.Range("L83").Value = 33491.5
$ws.Range("M83").Value = -384619258
$ws.Range("N83").Value = -43475.5
$ws.Range("H97").Value = 928.75
$ws.Range("I97").Value = 855.129
$ws.Range("J97").Value = 1182.3334
$ws.Range("K97").Value = 855.129
$ws.Range("L97").Value = 1182.3334
$ws.Range("M97").Value = -359.129
$ws.Range("N97").Value = -2174.3334

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4078.4783
$ws.Range("J7").Value = 6814.6665
$ws.Range("L7").Value = 6814.6665
$ws.Range("N7").Value = -7038.6665
$ws.Range("H40").Value = 4364.521
$ws.Range("I40").Value = 4312.5674
$ws.Range("K40").Value = 4312.5674
$ws.Range("M40").Value = -4176.5674
$ws.Range("H122").Value = 3371.2632
$ws.Range("I122").Value = 3214.8572
$ws.Range("J122").Value = 3809.2
$ws.Range("K122").Value = 9644.571599999999
$ws.Range("L122").Value = 11427.6
$ws.Range("M122").Value = -7194.571599999999
$ws.Range("N122").Value = -16327.6
$ws.Range("H126").Value = 4078.4783
$ws.Range("J126").Value = 6814.6665
$ws.Range("L126").Value = 20443.9995
$ws.Range("N126").Value = -25383.9995
$ws.Range("H132").Value = 2080.9443
$ws.Range("I132").Value = 2008.2858
$ws.Range("J132").Value = 2335.25
$ws.Range("K132").Value = 6024.857400000001
$ws.Range("L132").Value = 7005.75
$ws.Range("M132").Value = -3494.857400000001
$ws.Range("N132").Value = -12065.75
$ws.Range("H136").Value = 23080662
$ws.Range("I136").Value = 10207303
$ws.Range("K136").Value = 30621909
$ws.Range("M136").Value = -30619359

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 6705.8535
$ws.Range("I132").Value = 5148.3237
$ws.Range("J132").Value = 14271
$ws.Range("K132").Value = 15444.9711
$ws.Range("L132").Value = 42813
$ws.Range("M132").Value = -12914.9711
$ws.Range("N132").Value = -47873
$ws.Range("H136").Value = 8932468
$ws.Range("I136").Value = 11905824
$ws.Range("K136").Value = 35717472
$ws.Range("M136").Value = -35714922
